# Update "想去人数" (attendance) counts for two events, on both the
# "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1237
$wsExhibit.Range("F4").Value = 2729
$wsExhibit.Range("F5").Value = 243

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1237
$wsAll.Range("F6").Value = 2729
$wsAll.Range("F8").Value = 243
